# Updated cryptos list on Mon Aug 14 11:23:10 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.349.78'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +0.05%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.845.48'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -0.04%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9981'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '239.96'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6274'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.13%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9992'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07486'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -1.39%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2897'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.10%  '
$ws.Range('E10').Value = '  -0.94%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07731'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.00%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.845.72'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -2.25%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6798'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.24%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001049'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.83%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '82.02'
$ws.Range('D16').ClearFormats()
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.179'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.90%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '29.375.86'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.04%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '228.84'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.54%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.33'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.04%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9989'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.08%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.497'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.41%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9992'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '158.43'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.23%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.428'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.06%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1369'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.86%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.51'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.69%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06489'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +16.05%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.412'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -1.52%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.479'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +1.41%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.092'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.21%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.093'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.73%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.829'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.00%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.141'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.62%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6953'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.11%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.580'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.02%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.262.09'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +2.98%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.834'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +4.23%  '
$ws.Range('E39').Value = '  +2.00%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.763'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +6.48%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9191'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +2.58%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9984'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.008.64'
$ws.Range('D43').ClearFormats()
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.24'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.07%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '66.19'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +1.15%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.080'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -1.68%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.724'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +2.15%  '
$ws.Range('B48').Value = 'Algorand'
$ws.Range('C48').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1162'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +1.97%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.983'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.69%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.3957'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.71%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05695'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.02%  '
